$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 640.2
$ws.Range("C2").Value = 622.3
$ws.Range("D2").Value = 628.2
$ws.Range("E2").Value = 628.7
$ws.Range("F2").Value = 14
$ws.Range("G2").Value = 635

# Row 3
$ws.Range("B3").Value = 3038
$ws.Range("C3").Value = 2981.3
$ws.Range("D3").Value = 3022
$ws.Range("E3").Value = 3020.15
$ws.Range("F3").Value = 18
$ws.Range("G3").Value = 3034.95

# Row 4
$ws.Range("B4").Value = 500.5
$ws.Range("C4").Value = 485
$ws.Range("D4").Value = 493.1
$ws.Range("E4").Value = 492.35
$ws.Range("F4").Value = 25
$ws.Range("G4").Value = 498.6

# Row 5
$ws.Range("B5").Value = 1772
$ws.Range("C5").Value = 1715.9
$ws.Range("D5").Value = 1758.3
$ws.Range("E5").Value = 1755.65
$ws.Range("F5").Value = 55
$ws.Range("G5").Value = 1718.55

# Row 6
$ws.Range("B6").Value = 7158
$ws.Range("C6").Value = 6893.2
$ws.Range("D6").Value = 7066.65
$ws.Range("E6").Value = 7063.55
$ws.Range("F6").Value = 26
$ws.Range("G6").Value = 6895.05

# Row 7
$ws.Range("B7").Value = 196.5
$ws.Range("C7").Value = 191.64
$ws.Range("D7").Value = 194.22
$ws.Range("E7").Value = 193.93
$ws.Range("F7").Value = 212
$ws.Range("G7").Value = 196.27

# Row 8
$ws.Range("B8").Value = 250.7
$ws.Range("C8").Value = 247.65
$ws.Range("D8").Value = 250.1
$ws.Range("E8").Value = 249.8
$ws.Range("F8").Value = 99
$ws.Range("G8").Value = 249.65

# Row 9
$ws.Range("B9").Value = 530
$ws.Range("C9").Value = 516.1
$ws.Range("D9").Value = 529.15
$ws.Range("E9").Value = 527.55
$ws.Range("F9").Value = 91
$ws.Range("G9").Value = 525.35

# Row 10
$ws.Range("B10").Value = 835.35
$ws.Range("C10").Value = 818.5
$ws.Range("D10").Value = 835.2
$ws.Range("E10").Value = 831.9
$ws.Range("F10").Value = 43
$ws.Range("G10").Value = 834.05

# Row 11
$ws.Range("B11").Value = 4994.4
$ws.Range("C11").Value = 4880
$ws.Range("D11").Value = 4902.65
$ws.Range("E11").Value = 4915.9
$ws.Range("F11").Value = 5
$ws.Range("G11").Value = 4945

# Row 12
$ws.Range("B12").Value = 196.88
$ws.Range("C12").Value = 194.5
$ws.Range("D12").Value = 195.99
$ws.Range("E12").Value = 195.92
$ws.Range("F12").Value = 55
$ws.Range("G12").Value = 195.28

# Row 13
$ws.Range("B13").Value = 1756
$ws.Range("C13").Value = 1707.55
$ws.Range("D13").Value = 1750.45
$ws.Range("E13").Value = 1751.85
$ws.Range("F13").Value = 51
$ws.Range("G13").Value = 1709.4

# Row 14
$ws.Range("B14").Value = 1644.4
$ws.Range("C14").Value = 1631.15
$ws.Range("D14").Value = 1639.95
$ws.Range("E14").Value = 1638.55
$ws.Range("F14").Value = 245
$ws.Range("G14").Value = 1644.15

# Row 15
$ws.Range("B15").Value = 702.35
$ws.Range("C15").Value = 691.4
$ws.Range("D15").Value = 700
$ws.Range("E15").Value = 700.5
$ws.Range("F15").Value = 57
$ws.Range("G15").Value = 695.8

# Row 16
$ws.Range("B16").Value = 1235
$ws.Range("C16").Value = 1218.9
$ws.Range("D16").Value = 1222.75
$ws.Range("E16").Value = 1221.9
$ws.Range("F16").Value = 108
$ws.Range("G16").Value = 1222.8

# Row 17
$ws.Range("B17").Value = 1423.05
$ws.Range("C17").Value = 1407.8
$ws.Range("D17").Value = 1415.5
$ws.Range("E17").Value = 1417.45
$ws.Range("F17").Value = 33
$ws.Range("G17").Value = 1413.4

# Row 18
$ws.Range("B18").Value = 1950.5
$ws.Range("C18").Value = 1920.9
$ws.Range("D18").Value = 1930.2
$ws.Range("E18").Value = 1933.35
$ws.Range("F18").Value = 72
$ws.Range("G18").Value = 1929.05

# Row 19
$ws.Range("B19").Value = 967
$ws.Range("C19").Value = 948.55
$ws.Range("D19").Value = 961
$ws.Range("E19").Value = 960.5
$ws.Range("F19").Value = 10
$ws.Range("G19").Value = 962.35

# Row 20
$ws.Range("B20").Value = 681.9
$ws.Range("C20").Value = 667.05
$ws.Range("D20").Value = 677.75
$ws.Range("E20").Value = 675.95
$ws.Range("F20").Value = 13
$ws.Range("G20").Value = 676.95

# Row 21
$ws.Range("B21").Value = 2801.65
$ws.Range("C21").Value = 2741.65
$ws.Range("D21").Value = 2763
$ws.Range("E21").Value = 2757.6
$ws.Range("F21").Value = 22
$ws.Range("G21").Value = 2791

# Row 22
$ws.Range("B22").Value = 317.8
$ws.Range("C22").Value = 308.9
$ws.Range("D22").Value = 314.7
$ws.Range("E22").Value = 313.4
$ws.Range("F22").Value = 45
$ws.Range("G22").Value = 315.75

# Row 23
$ws.Range("B23").Value = 412.6
$ws.Range("C23").Value = 403.75
$ws.Range("D23").Value = 411.1
$ws.Range("E23").Value = 409.9
$ws.Range("F23").Value = 190
$ws.Range("G23").Value = 410.5

# Row 24
$ws.Range("B24").Value = 3074
$ws.Range("C24").Value = 2988
$ws.Range("D24").Value = 3042.9
$ws.Range("E24").Value = 3041.85
$ws.Range("F24").Value = 191
$ws.Range("G24").Value = 2997.15

# Row 25
$ws.Range("B25").Value = 816
$ws.Range("C25").Value = 809.1
$ws.Range("D25").Value = 815.7
$ws.Range("E25").Value = 814.5
$ws.Range("F25").Value = 182
$ws.Range("G25").Value = 810.4

# Row 26
$ws.Range("B26").Value = 815.95
$ws.Range("C26").Value = 800.55
$ws.Range("D26").Value = 810.65
$ws.Range("E26").Value = 809.9
$ws.Range("F26").Value = 14
$ws.Range("G26").Value = 814.05

# Row 27
$ws.Range("B27").Value = 1076.55
$ws.Range("C27").Value = 1057.55
$ws.Range("D27").Value = 1072
$ws.Range("E27").Value = 1073.3
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 1073.85

# Row 28
$ws.Range("B28").Value = 1142
$ws.Range("C28").Value = 1066
$ws.Range("D28").Value = 1112.9
$ws.Range("E28").Value = 1121.65
$ws.Range("F28").Value = 406
$ws.Range("G28").Value = 1076.25

# Row 29
$ws.Range("B29").Value = 433.1
$ws.Range("C29").Value = 424.7
$ws.Range("D29").Value = 430.9
$ws.Range("E29").Value = 430.9
$ws.Range("F29").Value = 99
$ws.Range("G29").Value = 430.4

# Row 30
$ws.Range("B30").Value = 153.56
$ws.Range("C30").Value = 151.2
$ws.Range("D30").Value = 153.15
$ws.Range("E30").Value = 152.97
$ws.Range("F30").Value = 444
$ws.Range("G30").Value = 152.85

# Row 31
$ws.Range("B31").Value = 11256.7
$ws.Range("C31").Value = 11115.4
$ws.Range("D31").Value = 11210
$ws.Range("E31").Value = 11222
$ws.Range("F31").Value = 4
$ws.Range("G31").Value = 11141

